# Move the THEOMAR (account 004231509) balance row:
#  - remove it from its current location (row 58, value 296.63)
#  - re-insert it just before the MARCO (004454491) row (row 33),
#    with the updated balance 1096.63

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row first (it is further down the sheet, so deleting it
# does not affect the row numbering of the insertion point above it).
$ws.Rows.Item(58).Delete()

# Insert a new row above the MARCO / 004454491 row (currently row 33)
# and populate it with the THEOMAR data at its new balance.
$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "004231509"
$ws.Cells.Item(33, 2).Value = "THEOMAR"
$ws.Cells.Item(33, 3).Value = 1096.63
